# Edit script: reorders several same-day match rows (swap F:V content between
# paired rows while keeping the A-E index/metadata columns fixed) and appends
# three newly scraped match rows (119-121) at the bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Swap the F:V (match data) content between each pair of rows that share
#    the same match date/time (column E) but were re-ordered by the scraper.
# ---------------------------------------------------------------------------
$pairs = @(
    @(10,11),
    @(14,15),
    @(16,17),
    @(29,30),
    @(34,35),
    @(48,49),
    @(54,55),
    @(59,60),
    @(62,63),
    @(70,71),
    @(91,92),
    @(102,103),
    @(106,107),
    @(108,109)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("F$r1`:V$r1")
    $range2 = $ws.Range("F$r2`:V$r2")

    $vals1 = $range1.Value2
    $vals2 = $range2.Value2

    $range1.Value2 = $vals2
    $range2.Value2 = $vals1
}

# ---------------------------------------------------------------------------
# 2) Append three new match rows (119, 120, 121) after the previous last row
#    (118), copying the formatting of row 118 first so the styles (bold/
#    bordered index column, date-time number format, etc.) line up.
# ---------------------------------------------------------------------------
$ws.Range("A118:V118").Copy($ws.Range("A119:V121"))

$newRows = @(
    @{
        Row = 119; Idx = 118; E = 45242.42708333334
        F = "Dukla Prague B"; G = 2; H = "Plzen B"; I = 2
        J = 2.32; K = "12/11/2023 00:12"; L = 3.31; M = "12/11/2023 10:00"
        N = 3.38; O = "12/11/2023 00:12"; P = 4.1;  Q = "12/11/2023 09:56"
        R = 2.68; S = "12/11/2023 00:12"; T = 1.85; U = "12/11/2023 10:00"
        V = "https://www.betexplorer.com/football/czech-republic/cfl-group-a/dukla-prague-plzen/AJ1ZO4g9/"
    },
    @{
        Row = 120; Idx = 119; E = 45242.42708333334
        F = "Slavia Prague B"; G = 2; H = "Hostoun"; I = 1
        J = 1.41; K = "12/11/2023 00:12"; L = 1.28; M = "12/11/2023 09:48"
        N = 4.7;  O = "12/11/2023 00:12"; P = 5.43; Q = "12/11/2023 09:51"
        R = 5.58; S = "12/11/2023 00:12"; T = 8;    U = "12/11/2023 09:48"
        V = "https://www.betexplorer.com/football/czech-republic/cfl-group-a/slavia-prague-hostoun/4MQZ33Hq/"
    },
    @{
        Row = 121; Idx = 120; E = 45242.58333333334
        F = "Povltavska FA"; G = 1; H = "Bohemians 1905 B"; I = 2
        J = 1.85; K = "12/11/2023 12:17"; L = 2.3;  M = "12/11/2023 13:57"
        N = 3.7;  O = "12/11/2023 12:17"; P = 3.88; Q = "12/11/2023 13:57"
        R = 3.66; S = "12/11/2023 12:17"; T = 2.55; U = "12/11/2023 13:57"
        V = "https://www.betexplorer.com/football/czech-republic/cfl-group-a/povltavska-fa-bohemians-1905/McLv3NWk/"
    }
)

foreach ($rd in $newRows) {
    $r = $rd.Row
    $ws.Range("A$r").Value = $rd.Idx
    $ws.Range("B$r").Value = "czech-republic"
    $ws.Range("C$r").Value = "cfl-group-a"
    $ws.Range("D$r").Value = "2023-2024"
    $ws.Range("E$r").Value = $rd.E
    $ws.Range("F$r").Value = $rd.F
    $ws.Range("G$r").Value = $rd.G
    $ws.Range("H$r").Value = $rd.H
    $ws.Range("I$r").Value = $rd.I
    $ws.Range("J$r").Value = $rd.J
    $ws.Range("K$r").Value = $rd.K
    $ws.Range("L$r").Value = $rd.L
    $ws.Range("M$r").Value = $rd.M
    $ws.Range("N$r").Value = $rd.N
    $ws.Range("O$r").Value = $rd.O
    $ws.Range("P$r").Value = $rd.P
    $ws.Range("Q$r").Value = $rd.Q
    $ws.Range("R$r").Value = $rd.R
    $ws.Range("S$r").Value = $rd.S
    $ws.Range("T$r").Value = $rd.T
    $ws.Range("U$r").Value = $rd.U
    $ws.Range("V$r").Value = $rd.V
}

Write-Host "done"
